$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D cells whose new price text looks like a plain number (e.g. "0.450", "1.00", "26.00").
# Excel auto-converts such text to a numeric value on assignment, but the source data
# keeps these as literal text, so pre-format those cells as Text before assigning them.
$textPriceRows = @(4,5,6,9,10,12,15,19,20,21,22,23,24,25,30,31,33,34,36,37,38,40,41,43,44,46,50,51)
foreach ($r in $textPriceRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

# Apply the updated price (column D) and 1h volume change (column E) values
$ws.Range("D2").Value = "57.539.06"
$ws.Range("E2").Value = "  -4.26%  "
$ws.Range("D3").Value = "3.142.73"
$ws.Range("E3").Value = "  -5.06%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "525.67"
$ws.Range("E5").Value = "  -6.12%  "
$ws.Range("D6").Value = "132.92"
$ws.Range("E6").Value = "  -7.94%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").Value = "3.139.63"
$ws.Range("E8").Value = "  -5.12%  "
$ws.Range("D9").Value = "0.450"
$ws.Range("E9").Value = "  -7.02%  "
$ws.Range("D10").Value = "7.21"
$ws.Range("E10").Value = "  -7.70%  "
$ws.Range("E11").Value = "  -7.94%  "
$ws.Range("D12").Value = "0.390"
$ws.Range("E12").Value = "  -4.54%  "
$ws.Range("D13").Value = "3.679.16"
$ws.Range("E13").Value = "  -5.70%  "
$ws.Range("E14").Value = "  -1.91%  "
$ws.Range("D15").Value = "25.88"
$ws.Range("E15").Value = "  -5.47%  "
$ws.Range("D16").Value = "3.142.28"
$ws.Range("E16").Value = "  -5.95%  "
$ws.Range("D17").Value = "57.532.20"
$ws.Range("E17").Value = "  -4.54%  "
$ws.Range("E18").Value = "  -8.25%  "
$ws.Range("D19").Value = "5.79"
$ws.Range("E19").Value = "  -6.31%  "
$ws.Range("D20").Value = "13.08"
$ws.Range("E20").Value = "  -8.65%  "
$ws.Range("D21").Value = "8.03"
$ws.Range("E21").Value = "  -8.36%  "
$ws.Range("D22").Value = "347.05"
$ws.Range("E22").Value = "  -7.10%  "
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "69.39"
$ws.Range("E24").Value = "  -6.54%  "
$ws.Range("D25").Value = "0.509"
$ws.Range("E25").Value = "  -7.60%  "
$ws.Range("D26").Value = "3.275.56"
$ws.Range("E26").Value = "  -6.04%  "
$ws.Range("D27").Value = "0.0₃0958"
$ws.Range("E27").Value = "  -9.26%  "
$ws.Range("E28").Value = "  -3.63%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").Value = "6.81"
$ws.Range("E30").Value = "  -5.67%  "
$ws.Range("D31").Value = "0.998"
$ws.Range("E31").Value = "  -0.27%  "
$ws.Range("E32").Value = "  -8.98%  "
$ws.Range("D33").Value = "6.93"
$ws.Range("E33").Value = "  -8.77%  "
$ws.Range("D34").Value = "21.63"
$ws.Range("E34").Value = "  -4.31%  "
$ws.Range("E35").Value = "  -4.99%  "
$ws.Range("D36").Value = "4.94"
$ws.Range("E36").Value = "  -5.25%  "
$ws.Range("D37").Value = "158.90"
$ws.Range("E37").Value = "  -4.77%  "
$ws.Range("D38").Value = "6.23"
$ws.Range("E38").Value = "  -7.62%  "
$ws.Range("E39").Value = "  -8.15%  "
$ws.Range("D40").Value = "26.00"
$ws.Range("E40").Value = "  -5.39%  "
$ws.Range("D41").Value = "0.0695"
$ws.Range("E41").Value = "  -5.44%  "
$ws.Range("D42").Value = "3.167.46"
$ws.Range("E42").Value = "  -5.99%  "
$ws.Range("D43").Value = "40.33"
$ws.Range("E43").Value = "  -3.95%  "
$ws.Range("D44").Value = "0.691"
$ws.Range("E44").Value = "  -8.08%  "
$ws.Range("E45").Value = "  -4.08%  "
$ws.Range("D46").Value = "3.95"
$ws.Range("E46").Value = "  -6.50%  "
$ws.Range("E47").Value = "  -0.33%  "
$ws.Range("E48").Value = "  -7.88%  "
$ws.Range("D49").Value = "2.261.68"
$ws.Range("E49").Value = "  -7.00%  "
$ws.Range("D50").Value = "6.23"
$ws.Range("E50").Value = "  -5.74%  "
$ws.Range("D51").Value = "20.46"
$ws.Range("E51").Value = "  -6.38%  "

# Restore the default (unformatted) cell style now that the text values are locked in,
# matching the rest of the sheet which does not carry an explicit number format.
foreach ($r in $textPriceRows) {
    $ws.Range("D$r").ClearFormats()
}

Write-Host "Updated cryptos list values"
